$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Run 50" column (AZ). This shifts the "Mean" column (old BA)
# left into AZ, matching the target layout (A1:AZ14 instead of A1:BA14).
$ws.Columns("AZ").Delete()

# Rename the header of column A from "Gen" to "MaxFES".
$ws.Range("A1").Value = "MaxFES"

# Update column A (row labels) from generation counts to MaxFES fractions.
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Update the recalculated "Mean" values now living in column AZ.
$ws.Range("AZ2").Value = 507344966.1302503
$ws.Range("AZ3").Value = 233134326.0406851
$ws.Range("AZ4").Value = 27714021.41196866
$ws.Range("AZ5").Value = 950100.04438669
$ws.Range("AZ6").Value = 312768.51985661
$ws.Range("AZ7").Value = 141938.09405402
$ws.Range("AZ8").Value = 79350.36870568
$ws.Range("AZ9").Value = 45688.2714633
$ws.Range("AZ10").Value = 29498.6900939
$ws.Range("AZ11").Value = 19533.4238549
$ws.Range("AZ12").Value = 14870.31868651
$ws.Range("AZ13").Value = 11233.11763967
$ws.Range("AZ14").Value = 8844.987848160001
